$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.243.22"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "2.746.04"
$ws.Range("E3").Value = "  +3.56%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'604.97"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'166.74"
$ws.Range("E6").Value = "  +5.20%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D9").Value = "2.743.24"
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").Value = "'0.143"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'0.365"
$ws.Range("E11").Value = "  +3.71%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "'28.81"
$ws.Range("E14").Value = "  +2.45%  "
$ws.Range("D15").Value = "3.241.70"
$ws.Range("E15").Value = "  +3.51%  "
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "69.147.66"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "2.732.55"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("D19").Value = "'11.97"
$ws.Range("E19").Value = "  +4.80%  "
$ws.Range("E20").Value = "  +5.35%  "
$ws.Range("D21").Value = "'369.36"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("D22").Value = "'4.58"
$ws.Range("E22").Value = "  +3.29%  "
$ws.Range("D23").Value = "'5.00"
$ws.Range("E23").Value = "  +3.69%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("D25").Value = "'74.27"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("D28").Value = "2.865.24"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").Value = "'603.73"
$ws.Range("E30").Value = "  +8.00%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("D32").Value = "'8.37"
$ws.Range("E32").Value = "  +4.23%  "
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("D36").Value = "'1.65"
$ws.Range("E36").Value = "  +4.90%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'163.25"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").Value = "'20.18"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").Value = "'0.384"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").Value = "'5.54"
$ws.Range("E42").Value = "  +3.44%  "
$ws.Range("D43").Value = "'2.73"
$ws.Range("E43").Value = "  +4.41%  "
$ws.Range("D44").Value = "'18.05"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").Value = "0.0₆0319"
$ws.Range("E45").Value = "  -4.07%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'159.33"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("E48").Value = "  +5.60%  "
$ws.Range("E49").Value = "  +6.56%  "
$ws.Range("D50").Value = "'0.614"
$ws.Range("E50").Value = "  +8.17%  "
$ws.Range("D51").Value = "'22.29"
$ws.Range("E51").Value = "  +0.50%  "
